$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1576.9474
$ws.Range("I40").Value = 1436.2
$ws.Range("K40").Value = 1436.2
$ws.Range("M40").Value = -1261.2
$ws.Range("H62").Value = 1999
$ws.Range("I62").Value = 1641
$ws.Range("J62").Value = 2161.7273
$ws.Range("K62").Value = 1641
$ws.Range("L62").Value = 2161.7273
$ws.Range("M62").Value = -1017
$ws.Range("N62").Value = -3409.7273
$ws.Range("H65").Value = 1999
$ws.Range("I65").Value = 1641
$ws.Range("J65").Value = 2161.7273
$ws.Range("K65").Value = 8205
$ws.Range("L65").Value = 10808.6365
$ws.Range("M65").Value = -5085
$ws.Range("N65").Value = -17048.6365
$ws.Range("H100").Value = 2470.5186
$ws.Range("I100").Value = 1810.9
$ws.Range("J100").Value = 2858.5293
$ws.Range("K100").Value = 1810.9
$ws.Range("L100").Value = 2858.5293
$ws.Range("M100").Value = -1269.9
$ws.Range("N100").Value = -3940.5293
$ws.Range("H101").Value = 2099.6428
$ws.Range("J101").Value = 1221.25
$ws.Range("L101").Value = 3663.75
$ws.Range("N101").Value = -6907.75
$ws.Range("H129").Value = 954.46155
$ws.Range("J129").Value = 1002.0175
$ws.Range("L129").Value = 3006.0525
$ws.Range("N129").Value = -13006.0525
$ws.Range("H132").Value = 2145.4082
$ws.Range("I132").Value = 1538.5483
$ws.Range("J132").Value = 3190.5557
$ws.Range("K132").Value = 4615.644899999999
$ws.Range("L132").Value = 9571.667099999999
$ws.Range("M132").Value = -2085.644899999999
$ws.Range("N132").Value = -14631.6671
$ws.Range("H137").Value = 1033.2808
$ws.Range("I137").Value = 824.0278
$ws.Range("J137").Value = 1392
$ws.Range("K137").Value = 2472.0834
$ws.Range("L137").Value = 4176
$ws.Range("M137").Value = 77.91660000000002
$ws.Range("N137").Value = -9276
$ws.Range("H138").Value = 3825.0422
$ws.Range("I138").Value = 1906.5652
$ws.Range("J138").Value = 7355.04
$ws.Range("K138").Value = 5719.6956
$ws.Range("L138").Value = 22065.12
$ws.Range("M138").Value = -579.6956
$ws.Range("N138").Value = -32345.12
$ws.Range("H141").Value = 3486.9648
$ws.Range("I141").Value = 1680.2264
$ws.Range("J141").Value = 27426.25
$ws.Range("K141").Value = 5040.6792
$ws.Range("L141").Value = 82278.75
$ws.Range("M141").Value = 139.3207999999995
$ws.Range("N141").Value = -92638.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1696.2354
$ws.Range("I61").Value = 1205.5
$ws.Range("J61").Value = 3986.3333
$ws.Range("K61").Value = 1205.5
$ws.Range("L61").Value = 3986.3333
$ws.Range("M61").Value = -993.5
$ws.Range("N61").Value = -4410.3333
$ws.Range("H74").Value = 1182.3158
$ws.Range("I74").Value = 1133.375
$ws.Range("J74").Value = 1443.3334
$ws.Range("K74").Value = 1133.375
$ws.Range("L74").Value = 1443.3334
$ws.Range("M74").Value = -259.375
$ws.Range("N74").Value = -3191.3334
$ws.Range("H77").Value = 1182.3158
$ws.Range("I77").Value = 1133.375
$ws.Range("J77").Value = 1443.3334
$ws.Range("K77").Value = 5666.875
$ws.Range("L77").Value = 7216.666999999999
$ws.Range("M77").Value = -1298.875
$ws.Range("N77").Value = -15952.667
$ws.Range("H102").Value = 73294.21000000001
$ws.Range("I102").Value = 1975.9
$ws.Range("J102").Value = 251590
$ws.Range("K102").Value = 1975.9
$ws.Range("L102").Value = 251590
$ws.Range("M102").Value = -353.9000000000001
$ws.Range("N102").Value = -254834
$ws.Range("H132").Value = 2085.361
$ws.Range("I132").Value = 1303.5714
$ws.Range("J132").Value = 3179.8667
$ws.Range("K132").Value = 3910.7142
$ws.Range("L132").Value = 9539.6001
$ws.Range("M132").Value = -1380.7142
$ws.Range("N132").Value = -14599.6001
$ws.Range("H136").Value = 1696.2354
$ws.Range("I136").Value = 1205.5
$ws.Range("J136").Value = 3986.3333
$ws.Range("K136").Value = 3616.5
$ws.Range("L136").Value = 11958.9999
$ws.Range("M136").Value = -1066.5
$ws.Range("N136").Value = -17058.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 25250
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1887
$ws.Range("H94").Value = 127737.375
$ws.Range("I94").Value = 1369.6666
$ws.Range("J94").Value = 203558
$ws.Range("K94").Value = 1369.6666
$ws.Range("L94").Value = 203558
$ws.Range("M94").Value = -918.6666
$ws.Range("N94").Value = -204460
$ws.Range("H99").Value = 1585.7142
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 2100
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 2100
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -5096
$ws.Range("H134").Value = 1936.2954
$ws.Range("I134").Value = 1620.5897
$ws.Range("J134").Value = 4398.8
$ws.Range("K134").Value = 4861.7691
$ws.Range("L134").Value = 13196.4
$ws.Range("M134").Value = -2326.7691
$ws.Range("N134").Value = -18266.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1795.5
$ws.Range("I31").Value = 1268.7778
$ws.Range("J31").Value = 3375.6667
$ws.Range("K31").Value = 1268.7778
$ws.Range("L31").Value = 3375.6667
$ws.Range("M31").Value = -973.7778000000001
$ws.Range("N31").Value = -3965.6667
$ws.Range("H34").Value = 1795.5
$ws.Range("I34").Value = 1268.7778
$ws.Range("J34").Value = 3375.6667
$ws.Range("K34").Value = 1268.7778
$ws.Range("L34").Value = 3375.6667
$ws.Range("M34").Value = -1066.7778
$ws.Range("N34").Value = -3779.6667
$ws.Range("H58").Value = 843049.1
$ws.Range("I58").Value = 1684760.8
$ws.Range("J58").Value = 1337.4546
$ws.Range("K58").Value = 1684760.8
$ws.Range("L58").Value = 1337.4546
$ws.Range("M58").Value = -1684557.8
$ws.Range("N58").Value = -1743.4546
$ws.Range("H132").Value = 266391.72
$ws.Range("I132").Value = 347468.47
$ws.Range("J132").Value = 2892.3333
$ws.Range("K132").Value = 1042405.41
$ws.Range("L132").Value = 8676.999899999999
$ws.Range("M132").Value = -1039875.41
$ws.Range("N132").Value = -13736.9999
$ws.Range("H134").Value = 1793.6052
$ws.Range("I134").Value = 1365.1936
$ws.Range("J134").Value = 3690.8572
$ws.Range("K134").Value = 4095.5808
$ws.Range("L134").Value = 11072.5716
$ws.Range("M134").Value = -1560.5808
$ws.Range("N134").Value = -16142.5716
$ws.Range("H136").Value = 843049.1
$ws.Range("I136").Value = 1684760.8
$ws.Range("J136").Value = 1337.4546
$ws.Range("K136").Value = 5054282.4
$ws.Range("L136").Value = 4012.3638
$ws.Range("M136").Value = -5051732.4
$ws.Range("N136").Value = -9112.363799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 900
$ws.Range("M4").Value = -788
$ws.Range("H131").Value = 13001612
$ws.Range("J131").Value = 14720764
$ws.Range("L131").Value = 44162292
$ws.Range("N131").Value = -44172372
$ws.Range("H132").Value = 2071.8333
$ws.Range("I132").Value = 680.5
$ws.Range("K132").Value = 6124.5
$ws.Range("M132").Value = -3594.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 24999.77
$ws.Range("J5").Value = 24999.75
$ws.Range("L5").Value = 24999.75
$ws.Range("N5").Value = -25223.75
$ws.Range("H132").Value = 1270.125
$ws.Range("I132").Value = 749.13043
$ws.Range("J132").Value = 2601.5557
$ws.Range("K132").Value = 2247.39129
$ws.Range("L132").Value = 7804.6671
$ws.Range("M132").Value = 282.60871
$ws.Range("N132").Value = -12864.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 27185.125
$ws.Range("I61").Value = 35313.5
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 35313.5
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -35111.5
$ws.Range("N61").Value = -3204
$ws.Range("H93").Value = 827.7646999999999
$ws.Range("I93").Value = 539.3333
$ws.Range("K93").Value = 539.3333
$ws.Range("M93").Value = 708.6667
$ws.Range("H100").Value = 12399.8
$ws.Range("I100").Value = 17999.666
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 17999.666
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -17458.666
$ws.Range("N100").Value = -5082
$ws.Range("H113").Value = 27185.125
$ws.Range("I113").Value = 35313.5
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 35313.5
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = -33143.5
$ws.Range("N113").Value = -7140
$ws.Range("H132").Value = 3532.5117
$ws.Range("I132").Value = 3137.4
$ws.Range("J132").Value = 4444.3076
$ws.Range("K132").Value = 9412.200000000001
$ws.Range("L132").Value = 13332.9228
$ws.Range("M132").Value = -6882.200000000001
$ws.Range("N132").Value = -18392.9228
$ws.Range("H136").Value = 3003.4058
$ws.Range("I136").Value = 3164.3962
$ws.Range("J136").Value = 2470.125
$ws.Range("K136").Value = 9493.188600000001
$ws.Range("L136").Value = 7410.375
$ws.Range("M136").Value = -6943.188600000001
$ws.Range("N136").Value = -12510.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 883.6842
$ws.Range("I132").Value = 618.5909
$ws.Range("J132").Value = 1780.9231
$ws.Range("K132").Value = 1855.7727
$ws.Range("L132").Value = 5342.7693
$ws.Range("M132").Value = 674.2273
$ws.Range("N132").Value = -10402.7693
$ws.Range("H136").Value = 1560.8536
$ws.Range("I136").Value = 1323.8235
$ws.Range("J136").Value = 2712.1428
$ws.Range("K136").Value = 3971.4705
$ws.Range("L136").Value = 8136.428400000001
$ws.Range("M136").Value = -1421.4705
$ws.Range("N136").Value = -13236.4284
